# Bug fix: import previsionnel tableau de bord
# Insert two new budget line rows ("Frais refacturables public" /
# "Frais refacturables privé") just above the "COÛTS MARGINAUX" section
# on the first sheet, pushing every row below down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows right before the old row 28 ("COÛTS MARGINAUX").
# Inserting at the same index twice pushes both new rows above it.
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# Fill in the new line labels.
$ws.Range("A28").Value = "Frais refacturables public"
$ws.Range("A29").Value = "Frais refacturables privé"

# Match the row height used by the other detail rows.
$ws.Rows.Item(28).RowHeight = 15
$ws.Rows.Item(29).RowHeight = 15

# Restore the view state (scroll position / active selection).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("C29").Select()
